# Install base iNZight packages
# - rename the existing sheet to "list-column"
# - add a new "two-row-header" sheet (Clippy data with a two-row header)
# - restore each sheet's selection so "two-row-header" ends up the active tab

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list-column"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "two-row-header"

# Row 1 : variable names
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "species"
$ws2.Range("C1").Value = "death"
$ws2.Range("D1").Value = "weight"

# Row 2 : variable annotations (order matches the original authoring order)
$ws2.Range("D2").Value = "(in grams)"
$ws2.Range("B2").Value = "(office supply type)"
$ws2.Range("A2").Value = "(at birth)"
$ws2.Range("C2").Value = "(date is approximate)"

# Row 3 : the data itself
$ws2.Range("A3").Value = "Clippy"
$ws2.Range("B3").Value = "paperclip"
# Reuse the date cell's existing style/value from the original sheet
$ws1.Range("B4").Copy($ws2.Range("C3"))
$ws2.Range("D3").Value = 0.9

# Restore the selections recorded in each sheet (set sheet1's first so that
# sheet2 - set last - ends up as the active/selected tab)
$ws1.Range("A2:A5").Select() | Out-Null
$ws2.Range("A1:D1").Select() | Out-Null
